$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.124600112295781
$ws.Range("C2").Value = 0.03465183341640454
$ws.Range("D2").Value = 0.003026854329077899
$ws.Range("E2").Value = 0.06674794700206732
$ws.Range("F2").Value = 4.687133108435745
$ws.Range("J2").Value = 0.1815652751015673
$ws.Range("K2").Value = 1.585180651292717
$ws.Range("L2").Value = 0.2706623424019625
$ws.Range("M2").Value = 0.4434864055175893
$ws.Range("N2").Value = 4.954385763929437
$ws.Range("B3").Value = 2.095091171269985
$ws.Range("C3").Value = 0.03053409699843712
$ws.Range("D3").Value = 0.003034667849541428
$ws.Range("E3").Value = 0.0669464004178808
$ws.Range("F3").Value = 4.678831105748003
$ws.Range("J3").Value = 0.1820829666190562
$ws.Range("K3").Value = 1.553133237959116
$ws.Range("L3").Value = 0.2701311971988289
$ws.Range("M3").Value = 0.4391798820445239
$ws.Range("N3").Value = 4.962651873587788
$ws.Range("B4").Value = 2.078074266702032
$ws.Range("C4").Value = 0.02801690480482932
$ws.Range("D4").Value = 0.003041697496320772
$ws.Range("E4").Value = 0.06708069524934768
$ws.Range("F4").Value = 4.675267130709756
$ws.Range("J4").Value = 0.1824263713628191
$ws.Range("K4").Value = 1.534345845640416
$ws.Range("L4").Value = 0.2699131463232121
$ws.Range("M4").Value = 0.4367519697316027
$ws.Range("N4").Value = 4.968628249207086
$ws.Range("B5").Value = 2.071417062371182
$ws.Range("C5").Value = 0.02699386853328178
$ws.Range("D5").Value = 0.003045128268037267
$ws.Range("E5").Value = 0.06713855970852567
$ws.Range("F5").Value = 4.674200608056253
$ws.Range("J5").Value = 0.1825727468247642
$ws.Range("K5").Value = 1.526913707004326
$ws.Range("L5").Value = 0.2698515205998717
$ws.Range("M5").Value = 0.4358170429631869
$ws.Range("N5").Value = 4.97129024268223
$ws.Range("B6").Value = 2.070328396972457
$ws.Range("C6").Value = 0.02682415803437266
$ws.Range("D6").Value = 0.003045732282797697
$ws.Range("E6").Value = 0.06714835786953444
$ws.Range("F6").Value = 4.674046822314338
$ws.Range("J6").Value = 0.1825974414140692
$ws.Range("K6").Value = 1.525693131235755
$ws.Range("L6").Value = 0.2698429340735728
$ws.Range("M6").Value = 0.4356650909847914
$ws.Range("N6").Value = 4.971745949734597
$ws.Range("B7").Value = 2.077983362062412
$ws.Range("C7").Value = 0.02800309675768631
$ws.Range("D7").Value = 0.00304174146638081
$ws.Range("E7").Value = 0.06708146291070705
$ws.Range("F7").Value = 4.675251184742052
$ws.Range("J7").Value = 0.182428319360044
$ws.Range("K7").Value = 1.534244706596866
$ws.Range("L7").Value = 0.2699122048823099
$ws.Range("M7").Value = 0.4367391403477434
$ws.Range("N7").Value = 4.968663232395954
$ws.Range("B8").Value = 2.114196927825617
$ws.Range("C8").Value = 0.03322969673263287
$ws.Range("D8").Value = 0.003029088012956649
$ws.Range("E8").Value = 0.06681379654078778
$ws.Range("F8").Value = 4.683952420569966
$ws.Range("J8").Value = 0.18173848207981
$ws.Range("K8").Value = 1.573946068750047
$ws.Range("L8").Value = 0.2704568051335983
$ws.Range("M8").Value = 0.4419566626578231
$ws.Range("N8").Value = 4.957048935011983
$ws.Range("B9").Value = 2.193947096612817
$ws.Range("C9").Value = 0.04357022541915967
$ws.Range("D9").Value = 0.003021782674968776
$ws.Range("E9").Value = 0.06638724160042742
$ws.Range("F9").Value = 4.713178617154227
$ws.Range("J9").Value = 0.1805877909149949
$ws.Range("K9").Value = 1.658862247266512
$ws.Range("L9").Value = 0.2723801896234903
$ws.Range("M9").Value = 0.4539023922264818
$ws.Range("N9").Value = 4.941422483043652
$ws.Range("B10").Value = 2.257866694958125
$ws.Range("C10").Value = 0.05122788684069235
$ws.Range("D10").Value = 0.003026831374852712
$ws.Range("E10").Value = 0.06613327722510931
$ws.Range("F10").Value = 4.742067719299669
$ws.Range("J10").Value = 0.1798648009065875
$ws.Range("K10").Value = 1.725566787454284
$ws.Range("L10").Value = 0.2743124820920073
$ws.Range("M10").Value = 0.4637227066817147
$ws.Range("N10").Value = 4.93430282913306
$ws.Range("B11").Value = 2.288103256788816
$ws.Range("C11").Value = 0.05472578712617349
$ws.Range("D11").Value = 0.003031338557825336
$ws.Range("E11").Value = 0.06603053666031045
$ws.Range("F11").Value = 4.756821572446995
$ws.Range("J11").Value = 0.1795623159598598
$ws.Range("K11").Value = 1.756852895071972
$ws.Range("L11").Value = 0.2753038214120735
$ws.Range("M11").Value = 0.4684167409202189
$ws.Range("N11").Value = 4.932011715581297
$ws.Range("B12").Value = 2.299719653437592
$ws.Range("C12").Value = 0.05605249330488959
$ws.Range("D12").Value = 0.003033359036668237
$ws.Range("E12").Value = 0.06599346164800135
$ws.Range("F12").Value = 4.762640239423661
$ws.Range("J12").Value = 0.1794515574387177
$ws.Range("K12").Value = 1.768835643942509
$ws.Range("L12").Value = 0.2756953232582902
$ws.Range("M12").Value = 0.4702268103814973
$ws.Range("N12").Value = 4.931280451467188
$ws.Range("B13").Value = 2.297210454450749
$ws.Range("C13").Value = 0.05576666785810858
$ws.Range("D13").Value = 0.003032910001953937
$ws.Range("E13").Value = 0.06600136513221511
$ws.Range("F13").Value = 4.761376781971109
$ws.Range("J13").Value = 0.1794752430648643
$ws.Range("K13").Value = 1.766248924478049
$ws.Range("L13").Value = 0.2756102909237512
$ws.Range("M13").Value = 0.4698355332598751
$ws.Range("N13").Value = 4.931431877923472
$ws.Range("B14").Value = 2.28905560976159
$ws.Range("C14").Value = 0.05483489307577827
$ws.Range("D14").Value = 0.003031498518937248
$ws.Range("E14").Value = 0.06602744983675546
$ws.Range("F14").Value = 4.757295634603963
$ws.Range("J14").Value = 0.1795531279775915
$ws.Range("K14").Value = 1.75783601017028
$ws.Range("L14").Value = 0.2753357080556995
$ws.Range("M14").Value = 0.4685650046134029
$ws.Range("N14").Value = 4.931948821413073
$ws.Range("B15").Value = 2.284082204215622
$ws.Range("C15").Value = 0.05426443287119298
$ws.Range("D15").Value = 0.003030674677759748
$ws.Range("E15").Value = 0.06604366562141628
$ws.Range("F15").Value = 4.754825985474525
$ws.Range("J15").Value = 0.1796013274462016
$ws.Range("K15").Value = 1.752700485840052
$ws.Range("L15").Value = 0.2751696137626496
$ws.Range("M15").Value = 0.4677910050555738
$ws.Range("N15").Value = 4.932283220370749
$ws.Range("B16").Value = 2.255913971357757
$ws.Range("C16").Value = 0.05099958535612359
$ws.Range("D16").Value = 0.003026580942293577
$ws.Range("E16").Value = 0.06614024818817388
$ws.Range("F16").Value = 4.741135949515012
$ws.Range("J16").Value = 0.1798850994129051
$ws.Range("K16").Value = 1.723541114289503
$ws.Range("L16").Value = 0.2742499518341859
$ws.Range("M16").Value = 0.4634204974345977
$ws.Range("N16").Value = 4.934471632974137
$ws.Range("B17").Value = 2.238930412144043
$ws.Range("C17").Value = 0.04900043417136146
$ws.Range("D17").Value = 0.003024632721944087
$ws.Range("E17").Value = 0.06620276802622627
$ws.Range("F17").Value = 4.733150362289336
$ws.Range("J17").Value = 0.1800659396659707
$ws.Range("K17").Value = 1.705893975698416
$ws.Range("L17").Value = 0.2737145038059552
$ws.Range("M17").Value = 0.4607973574289872
$ws.Range("N17").Value = 4.936056910807153
$ws.Range("B18").Value = 2.229271042079461
$ws.Range("C18").Value = 0.04785192642616209
$ws.Range("D18").Value = 0.003023720459290224
$ws.Range("E18").Value = 0.06623993205519074
$ws.Range("F18").Value = 4.728709003989692
$ws.Range("J18").Value = 0.1801724404020586
$ws.Range("K18").Value = 1.695832465619191
$ws.Range("L18").Value = 0.2734171042084412
$ws.Range("M18").Value = 0.4593099406316412
$ws.Range("N18").Value = 4.937057914321827
$ws.Range("B19").Value = 2.226019294536343
$ws.Range("C19").Value = 0.04746329174345476
$ws.Range("D19").Value = 0.003023447488502029
$ws.Range("E19").Value = 0.06625272225028311
$ws.Range("F19").Value = 4.727231301377387
$ws.Range("J19").Value = 0.1802089270756912
$ws.Range("K19").Value = 1.692441040157433
$ws.Range("L19").Value = 0.2733182277686339
$ws.Range("M19").Value = 0.4588099946198199
$ws.Range("N19").Value = 4.937412154290939
$ws.Range("B20").Value = 2.240727048755502
$ws.Range("C20").Value = 0.0492131069745767
$ws.Range("D20").Value = 0.003024818585610767
$ws.Range("E20").Value = 0.06619598810670535
$ws.Range("F20").Value = 4.733984737928139
$ws.Range("J20").Value = 0.1800464316815926
$ws.Range("K20").Value = 1.707763368498803
$ws.Range("L20").Value = 0.2737704089543271
$ws.Range("M20").Value = 0.4610743865949232
$ws.Range("N20").Value = 4.935878923628721
$ws.Range("B21").Value = 2.29144636923229
$ws.Range("C21").Value = 0.05510851972105968
$ws.Range("D21").Value = 0.003031904623958859
$ws.Range("E21").Value = 0.06601973851145448
$ws.Range("F21").Value = 4.758488078988023
$ws.Range("J21").Value = 0.1795301486303771
$ws.Range("K21").Value = 1.760303413351323
$ws.Range("L21").Value = 0.2754159230580839
$ws.Range("M21").Value = 0.4689373070211005
$ws.Range("N21").Value = 4.931793282100358
$ws.Range("B22").Value = 2.32556457349483
$ws.Range("C22").Value = 0.05897393911791937
$ws.Range("D22").Value = 0.003038362459903965
$ws.Range("E22").Value = 0.06591521553856428
$ws.Range("F22").Value = 4.775852926325797
$ws.Range("J22").Value = 0.179214790133063
$ws.Range("K22").Value = 1.795430437648605
$ws.Range("L22").Value = 0.2765852005556368
$ws.Range("M22").Value = 0.4742658235654815
$ws.Range("N22").Value = 4.92991775065758
$ws.Range("B23").Value = 2.307266345851474
$ws.Range("C23").Value = 0.05690973583598691
$ws.Range("D23").Value = 0.00303474998392339
$ws.Range("E23").Value = 0.06597002821857156
$ws.Range("F23").Value = 4.766461443798335
$ws.Range("J23").Value = 0.1793810879011701
$ws.Range("K23").Value = 1.776610302738959
$ws.Range("L23").Value = 0.2759525651523305
$ws.Range("M23").Value = 0.4714045617622205
$ws.Range("N23").Value = 4.930846023776724
$ws.Range("B24").Value = 2.239914463175921
$ws.Range("C24").Value = 0.04911695501792224
$ws.Range("D24").Value = 0.00302473390928526
$ws.Range("E24").Value = 0.06619904950601896
$ws.Range("F24").Value = 4.733607050475001
$ws.Range("J24").Value = 0.1800552433501998
$ws.Range("K24").Value = 1.706917954151777
$ws.Range("L24").Value = 0.2737451017245149
$ws.Range("M24").Value = 0.4609490773455818
$ws.Range("N24").Value = 4.935959112476127
$ws.Range("B25").Value = 2.171437214005039
$ws.Range("C25").Value = 0.0407625729324792
$ws.Range("D25").Value = 0.00302191177864497
$ws.Range("E25").Value = 0.06649216269014779
$ws.Range("F25").Value = 4.703969869793028
$ws.Range("J25").Value = 0.1808775294856453
$ws.Range("K25").Value = 1.635132979706299
$ws.Range("L25").Value = 0.2717685249062356
$ws.Range("M25").Value = 0.4504873723807137
$ws.Range("N25").Value = 4.944884215423258
